# QA Round 2: deep quality optimization - compliance, diversification, UX improvements

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Small copy edits on the "EmilyBellJourney" sheet
# -----------------------------------------------------------------
$journey = $wb.Worksheets.Item("EmilyBellJourney")
$journey.Range("B5").Value  = "hold on hold on"
$journey.Range("B8").Value  = "oh god 😩"
$journey.Range("B20").Value = "so what do you think? 😏"
$journey.Range("B22").Value = "gimme a minute 🔥"

# -----------------------------------------------------------------
# 2) Split the "cumcontrol" sheet into "cumcontrol1" + "cumcontrol2"
#    - Duplicate the existing sheet (keeps formatting/col widths/styles)
#    - Rename the original to "cumcontrol1" and refresh its copy text
#    - Rename the duplicate to "cumcontrol2" and give it fresh variants
# -----------------------------------------------------------------
$cumcontrol = $wb.Worksheets.Item("cumcontrol")
$cumcontrol.Copy($null, $cumcontrol)

$cumcontrol.Name = "cumcontrol1"
$cumcontrol2 = $wb.Worksheets.Item("cumcontrol (2)")
$cumcontrol2.Name = "cumcontrol2"

# cumcontrol1 — refreshed copy for the original rows
$cumcontrol1 = $wb.Worksheets.Item("cumcontrol1")
$cumcontrol1.Range("B2").Value = "if you finish before you see what I'm sending next you'll regret it"

$cumcontrol1.Range("B3").Value = "wait wait wait babe... I have one more thing for you before you finish 💋"
$cumcontrol1.Range("C3").Value = "DELAY. Send final PPV."

$cumcontrol1.Range("B4").Value = "I want to feel it at the same time... watch this first"
$cumcontrol1.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol1.Range("B5").Value = "okay NOW we can go together... open this 🔥"
$cumcontrol1.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol1.Range("B6").Value = "you better not be close already babe... I have more to show you 💋"

$cumcontrol1.Range("B7").Value = "not yet... I said not yet 😏"
$cumcontrol1.Range("C7").Value = "CONTROL. More PPVs to send. Create urgency to open next."

# cumcontrol2 — brand-new variant copy (row names stay delay2/delay1/sync2/sync1/edge2/edge1)
$cumcontrol2.Range("B2").Value = "hold on just a little longer babe, I promise this next one is worth it 💋"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("B3").Value = "don't you dare... not until you see what I just did"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("B4").Value = "let's do this together... but you have to open this first"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("B5").Value = "okay I'm ready now too... watch this with me 🔥"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("B6").Value = "patience... the best part hasn't even happened yet"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("B7").Value = "slow down babe... I'm not letting you off that easy 😏"
$cumcontrol2.Range("C7").Value = "CONTROL."

# Restore original active sheet/selection
$journey.Activate() | Out-Null
